$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing "sum" header cell (G1) to the new "Save" header cell (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
